$wb = $excel.ActiveWorkbook

# --- Sheet "1er Parcial" (sheet1) ---
$ws1 = $wb.Worksheets.Item("1er Parcial")
$ws1.Range("E18").Value = 22
$ws1.Range("F18").Value = 3
$ws1.Range("G18").Value = 88
$ws1.Range("H18").Value = 12
$ws1.Range("I18").Value = 7.8
$ws1.Range("J18").Value = 3
$ws1.Range("K18").Value = 12

# --- Sheet "2o Parcial" (sheet2) ---
$ws2 = $wb.Worksheets.Item("2o Parcial")
$ws2.Range("E18").Value = 22
$ws2.Range("F18").Value = 3
$ws2.Range("G18").Value = 88
$ws2.Range("H18").Value = 12
$ws2.Range("I18").Value = 7.8
$ws2.Range("J18").Value = 3
$ws2.Range("K18").Value = 12

# --- Sheet "3er Parcial" (sheet3) ---
$ws3 = $wb.Worksheets.Item("3er Parcial")
$ws3.Range("E18").Value = 22
$ws3.Range("F18").Value = 3
$ws3.Range("G18").Value = 88
$ws3.Range("H18").Value = 12
$ws3.Range("J18").Value = 3
$ws3.Range("K18").Value = 12
